$d = $word.ActiveDocument

$d.Content.Find.Execute("45÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷8=", 2) | Out-Null
$d.Content.Find.Execute("41÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷3=", 2) | Out-Null
$d.Content.Find.Execute("49÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷5=", 2) | Out-Null
$d.Content.Find.Execute("86÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷2=", 2) | Out-Null
$d.Content.Find.Execute("68÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷2=", 2) | Out-Null
$d.Content.Find.Execute("80÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷9=", 2) | Out-Null
$d.Content.Find.Execute("71÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷5=", 2) | Out-Null
$d.Content.Find.Execute("43÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷3=", 2) | Out-Null
$d.Content.Find.Execute("78÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷7=", 2) | Out-Null
$d.Content.Find.Execute("43÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷6=", 2) | Out-Null
$d.Content.Find.Execute("75÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷2=", 2) | Out-Null
$d.Content.Find.Execute("69÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷5=", 2) | Out-Null
$d.Content.Find.Execute("85÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷2=", 2) | Out-Null
$d.Content.Find.Execute("14÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷3=", 2) | Out-Null
$d.Content.Find.Execute("61÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷4=", 2) | Out-Null
$d.Content.Find.Execute("84÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷5=", 2) | Out-Null
$d.Content.Find.Execute("88÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷8=", 2) | Out-Null
$d.Content.Find.Execute("62÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷4=", 2) | Out-Null
$d.Content.Find.Execute("93÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷3=", 2) | Out-Null
$d.Content.Find.Execute("76÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷2=", 2) | Out-Null
$d.Content.Find.Execute("86÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷9=", 2) | Out-Null
$d.Content.Find.Execute("54÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷3=", 2) | Out-Null
$d.Content.Find.Execute("48÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷8=", 2) | Out-Null
$d.Content.Find.Execute("55÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷9=", 2) | Out-Null
$d.Content.Find.Execute("89÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷2=", 2) | Out-Null
